# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.728.65'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '2.604.65'
$ws.Range('E3').Value = '  +0.76%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'570.61"
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').Value = "'142.88"
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('D7').Value = "'0.996"
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.59%  '
$ws.Range('D9').Value = '2.631.63'
$ws.Range('E9').Value = '  +1.38%  '
$ws.Range('D10').Value = "'6.50"
$ws.Range('E10').Value = '  -2.88%  '
$ws.Range('E11').Value = '  +2.74%  '
$ws.Range('D13').Value = "'0.367"
$ws.Range('E13').Value = '  +6.51%  '
$ws.Range('D14').Value = '3.071.22'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').Value = '60.712.81'
$ws.Range('E15').Value = '  +2.34%  '
$ws.Range('D16').Value = "'23.58"
$ws.Range('E16').Value = '  +4.52%  '
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '2.619.65'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = "'4.68"
$ws.Range('E19').Value = '  +3.11%  '
$ws.Range('D20').Value = "'11.20"
$ws.Range('E20').Value = '  +9.14%  '
$ws.Range('D21').Value = "'347.59"
$ws.Range('E21').Value = '  +3.43%  '
$ws.Range('D22').Value = "'7.09"
$ws.Range('E22').Value = '  +13.98%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = "'0.519"
$ws.Range('E24').Value = '  +13.99%  '
$ws.Range('D25').Value = "'63.76"
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').Value = "'0.996"
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'0.160"
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').Value = "'7.72"
$ws.Range('E28').Value = '  +5.71%  '
$ws.Range('D29').Value = '0.0₃0795'
$ws.Range('E29').Value = '  +1.69%  '
$ws.Range('D30').Value = "'1.85"
$ws.Range('E30').Value = '  +9.55%  '
$ws.Range('E31').Value = '  -0.11%  '
$ws.Range('D32').Value = "'6.29"
$ws.Range('E32').Value = '  +3.23%  '
$ws.Range('D33').Value = "'161.62"
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('D35').Value = "'4.26"
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('E36').Value = '  +8.82%  '
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('D38').Value = "'1.59"
$ws.Range('E38').Value = '  +6.23%  '
$ws.Range('D39').Value = "'37.69"
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('E40').Value = '  -2.26%  '
$ws.Range('E41').Value = '  +3.35%  '
$ws.Range('D42').Value = "'297.06"
$ws.Range('E42').Value = '  +0.74%  '
$ws.Range('D43').Value = "'139.20"
$ws.Range('E43').Value = '  +10.89%  '
$ws.Range('D44').Value = "'0.995"
$ws.Range('E44').Value = '  -0.40%  '
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').Value = "'0.0552"
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('D48').Value = "'0.0241"
$ws.Range('E48').Value = '  +3.49%  '
$ws.Range('D49').Value = "'10.69"
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').Value = "'19.66"
$ws.Range('E50').Value = '  +5.73%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.039.32'
$ws.Range('E51').Value = '  +4.42%  '
